# Checklist_Especificacao_do_Projeto_WheyMarket.xlsx
# "Add files via upload" — fills in the "RESPONSÁVEL" (E) and
# "CLASSIFICAÇÃO" (F) columns on the Checklist sheet for every
# checklist item row (3-24).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checklist")

$responsavel = "Leonardo Klein`nEduardo Moura"

# Row -> CLASSIFICAÇÃO text (column F), mirrors the matching
# classification already shown in column H for that bucket.
$classificacoes = @{
    3  = "MEDIA ALTA"
    4  = "ALTA"
    5  = "MEDIA"
    6  = "EXTREMAMENTE ALTA"
    7  = "MEDIA ALTA"
    8  = "EXTREMAMENTE ALTA"
    9  = "ALTA"
    10 = "MEDIA ALTA"
    11 = "MEDIA ALTA"
    12 = "EXTREMAMENTE ALTA"
    13 = "MEDIA ALTA"
    14 = "EXTREMAMENTE ALTA"
    15 = "MEDIA ALTA"
    16 = "ALTA"
    17 = "MEDIA ALTA"
    18 = "BAIXA"
    19 = "BAIXA"
    20 = "BAIXA"
    21 = "EXTREMAMENTE BAIXA"
    22 = "EXTREMAMENTE BAIXA"
    23 = "EXTREMAMENTE BAIXA"
    24 = "EXTREMAMENTE BAIXA"
}

foreach ($row in 3..24) {
    $eCell = $ws.Range("E$row")
    $eCell.Value = $responsavel
    $eCell.WrapText = $true
    $eCell.HorizontalAlignment = -4108
    $eCell.VerticalAlignment = -4108

    $fCell = $ws.Range("F$row")
    $fCell.Value = $classificacoes[$row]
    $fCell.WrapText = $true
    $fCell.HorizontalAlignment = -4108
    $fCell.VerticalAlignment = -4108
}

# Mirror the saved selection state on the Checklist sheet.
[void]$ws.Range("H19").Select()

Write-Output "done"
